# Team-Meeting-Agenda.docx edit script
# 1) Adds a "Download Word Document" hyperlink paragraph right after the
#    top "Back to Home | Guide and Rubric" line.
# 2) Adds a "Tip for Note Taker" note (with a few bold runs) right after
#    the "Share notes with team after the meeting" bullet / right before
#    the horizontal-rule paragraph that follows the Note Taker section.

$d = $word.ActiveDocument

# Keep the output close to the canonical OOXML - the source document has
# no rsid churn, so don't let the host stamp any in.
$word.Options.StoreRSIDOnSave = $false

# ---------------------------------------------------------------------
# Part 1: "Download Word Document" link under the title-bar navigation
# ---------------------------------------------------------------------

$topNav = $d.Paragraphs(1)
$ip = $topNav.Range
$ip.Collapse(0)
$ip.InsertParagraphAfter()

$dlPara = $d.Paragraphs(2)
$dlPara.Style = "Body Text"

$dlRange = $dlPara.Range
$null = $d.Hyperlinks.Add($dlRange, "team-meeting-agenda.docx", $null, $null, "Download Word Document")

# Make the link text bold, matching the source run formatting.
$dlHyperlink = $dlPara.Range.Hyperlinks(1)
$dlHyperlink.Range.Bold = 1

# ---------------------------------------------------------------------
# Part 2: "Tip for Note Taker" paragraph
# ---------------------------------------------------------------------

$shareNotes = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "Share notes with team after the meeting`r") {
        $shareNotes = $para
        break
    }
}

$ip2 = $shareNotes.Range
$ip2.Collapse(0)
$ip2.InsertParagraphAfter()

$tipPara = $d.Paragraphs($i + 1)
$tipPara.Style = "Block Text"

$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$tipText = "Tip for Note Taker" + ": " + "1. Click" + " " + $openQuote + "Download Word Document" + $closeQuote + " " + "at the top of this page to download the template" + " " + "2. Fill in the notes as your team discusses (sharing your screen is recommended)" + " " + "3. Keep notes brief and simple" + " " + "4." + " " + "Alternative" + ": If you record the meeting via Zoom, you can download the transcription after the meeting and ask AI to format it according to this meeting note template"

$tipPara.Range.Text = $tipText

$tipBase = $tipPara.Range.Start

# "Tip for Note Taker" (bold)
$run1 = $d.Range($tipBase + 0, $tipBase + 18)
$run1.Bold = 1

# the curly-quoted "Download Word Document" (bold)
$run2 = $d.Range($tipBase + 29, $tipBase + 53)
$run2.Bold = 1

# "Alternative" (bold)
$run3 = $d.Range($tipBase + 218, $tipBase + 229)
$run3.Bold = 1

Write-Host "Edit complete."
Write-Host "Download paragraph text:" $dlPara.Range.Text
Write-Host "Tip paragraph text:" $tipPara.Range.Text
